$d = $word.ActiveDocument

function Find-BookmarkByName($doc, $name) {
    for ($i = 1; $i -le $doc.Bookmarks.Count; $i++) {
        $b = $doc.Bookmarks.Item($i)
        if ($b.Name -eq $name) {
            return $b
        }
    }
    return $null
}

function Rename-Bookmark($doc, $oldName, $newName) {
    $bm = Find-BookmarkByName $doc $oldName
    if ($bm -eq $null) {
        return
    }
    $target = $doc.Range($bm.Start, $bm.End)
    # Re-adding a bookmark under its existing name moves it (Word keeps
    # bookmark names unique) instead of creating a duplicate, so shove the
    # old name out of the way of the text it used to wrap ...
    $origin = $doc.Range(0, 0)
    $doc.Bookmarks.Add($oldName, $origin) | Out-Null
    # ... and drop the new name bookmark exactly where the old one lived.
    $doc.Bookmarks.Add($newName, $target) | Out-Null
}

# Fix the two typos in the bolded instruction text.
$d.Content.Find.Execute(
    "Select affected locations by typing or scrolling through the selection men:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Select affected locations by typing or scrolling through the selection menu:",
    2) | Out-Null

$d.Content.Find.Execute(
    "The application will display the summarized information by indicating the Affected Areas and the appropriate reporting forma:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The application will display the summarized information by indicating the Affected Areas and the appropriate reporting format:",
    2) | Out-Null

# The auto-generated bookmarks anchored on those headings are derived from
# the heading text, so they need to be regenerated to match the corrected
# wording.
Rename-Bookmark $d "X698896adfa04d4b2a71935444958772c5b46463" "Xc338b168b7cccd99ed8cc21bb1f6b9ec4bb776a"
Rename-Bookmark $d "Xb00c727a6577662d9aec1f5253063c4ca479014" "Xb6f9d67f48d1f26911168aa063102d9371c2f99"
